# Update the "fflush(stdout);" / "//fflush(stdout);" entries in the
# StrFind / StrReplace table to drop the trailing semicolon, i.e.
# "fflush(stdout);" -> "fflush(stdout" and "//fflush(stdout);" -> "//fflush(stdout"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B22").Value = "fflush(stdout"
$ws.Range("C22").Value = "//fflush(stdout"

# Mirror the author's navigation: the selection ended up on B22 after the edit.
$ws.Range("B22").Select()
